$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 147, shifting the existing rows 147:165 down to 148:166.
$ws.Rows("147:147").Insert()

# Populate the newly inserted row 147 with a fresh weekly data point
# (same market/category attributes, latest date, last week's price values).
$ws.Range("A147").Value = 4
$ws.Range("B147").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C147").Value = "Los Lagos"
$ws.Range("D147").Value = 44491
$ws.Range("E147").Value = 10
$ws.Range("F147").Value = 100112044
$ws.Range("G147").Value = "Perejil"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 180
$ws.Range("K147").Value = 5000
$ws.Range("L147").Value = 5000
$ws.Range("M147").Value = 5000
$ws.Range("N147").Value = "$/docena de atados (3 kilos)"
$ws.Range("O147").Value = "Región Metropolitana"
$ws.Range("P147").Value = 1667
$ws.Range("Q147").Value = 3
$ws.Range("R147").Value = "Hortaliza"
